$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Range("H17").Value = 45456316
$ws_ALC.Range("J17").Value = 45456316
$ws_ALC.Range("L17").Value = 136368948
$ws_ALC.Range("N17").Value = -136369284

$ws_ALC.Range("H34").Value = 11666.667
$ws_ALC.Range("I34").Value = 11666.667
$ws_ALC.Range("J34").Value = 0
$ws_ALC.Range("K34").Value = 11666.667
$ws_ALC.Range("L34").Value = 0
$ws_ALC.Range("M34").ClearContents()
$ws_ALC.Range("N34").Value = -11463.667

$ws_ALC.Range("H36").Value = 11666.667
$ws_ALC.Range("I36").Value = 11666.667
$ws_ALC.Range("J36").Value = 0
$ws_ALC.Range("K36").Value = 11666.667
$ws_ALC.Range("L36").Value = 0
$ws_ALC.Range("M36").ClearContents()
$ws_ALC.Range("N36").Value = -10951.667

$ws_ALC.Range("H43").Value = 1787
$ws_ALC.Range("I43").Value = 1370.4667
$ws_ALC.Range("J43").Value = 2828.3333
$ws_ALC.Range("K43").Value = 1370.4667
$ws_ALC.Range("L43").Value = 2828.3333
$ws_ALC.Range("M43").Value = -1301.4667
$ws_ALC.Range("N43").Value = -2966.3333

$ws_ALC.Range("H62").Value = 4420.6
$ws_ALC.Range("I62").Value = 4420.6
$ws_ALC.Range("K62").Value = 4420.6
$ws_ALC.Range("M62").Value = -3796.6

$ws_ALC.Range("H64").Value = 7016.722
$ws_ALC.Range("I64").Value = 5849.6
$ws_ALC.Range("K64").Value = 5849.6
$ws_ALC.Range("M64").Value = -5601.6

$ws_ALC.Range("H65").Value = 4420.6
$ws_ALC.Range("I65").Value = 4420.6
$ws_ALC.Range("K65").Value = 22103
$ws_ALC.Range("M65").Value = -18983

$ws_ALC.Range("H67").Value = 7016.722
$ws_ALC.Range("I67").Value = 5849.6
$ws_ALC.Range("K67").Value = 5849.6
$ws_ALC.Range("M67").Value = -4991.6

$ws_ALC.Range("H101").Value = 499.55554
$ws_ALC.Range("I101").Value = 554.8
$ws_ALC.Range("J101").Value = 430.5
$ws_ALC.Range("K101").Value = 1664.4
$ws_ALC.Range("L101").Value = 1291.5
$ws_ALC.Range("M101").Value = -42.39999999999986
$ws_ALC.Range("N101").Value = -4535.5

$ws_ALC.Range("H106").Value = 8243.5
$ws_ALC.Range("I106").Value = 8243.5
$ws_ALC.Range("K106").Value = 8243.5
$ws_ALC.Range("M106").Value = -7612.5

$ws_ALC.Range("H116").Value = 3399.5
$ws_ALC.Range("I116").Value = 3399.5
$ws_ALC.Range("K116").Value = 3399.5
$ws_ALC.Range("M116").Value = 42.5

$ws_ALC.Range("H132").Value = 22539.084
$ws_ALC.Range("I132").Value = 1733.6316
$ws_ALC.Range("J132").Value = 101599.8
$ws_ALC.Range("K132").Value = 5200.8948
$ws_ALC.Range("L132").Value = 304799.4
$ws_ALC.Range("M132").Value = -2670.8948
$ws_ALC.Range("N132").Value = -309859.4

$ws_ALC.Range("H138").Value = 2037.7028
$ws_ALC.Range("I138").Value = 1201.6666
$ws_ALC.Range("J138").Value = 4295
$ws_ALC.Range("K138").Value = 3604.9998
$ws_ALC.Range("L138").Value = 12885
$ws_ALC.Range("M138").Value = 1535.0002
$ws_ALC.Range("N138").Value = -23165

$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Range("H32").Value = 9456.052
$ws_ARM.Range("I32").Value = 7403.091
$ws_ARM.Range("J32").Value = 20747.334
$ws_ARM.Range("K32").Value = 7403.091
$ws_ARM.Range("L32").Value = 20747.334
$ws_ARM.Range("M32").Value = -7116.091
$ws_ARM.Range("N32").Value = -21321.334

$ws_ARM.Range("H61").Value = 5862.4736
$ws_ARM.Range("I61").Value = 4360.3335
$ws_ARM.Range("J61").Value = 7214.4
$ws_ARM.Range("K61").Value = 4360.3335
$ws_ARM.Range("L61").Value = 7214.4
$ws_ARM.Range("M61").Value = -4148.3335
$ws_ARM.Range("N61").Value = -7638.4

$ws_ARM.Range("H88").Value = 2015.75
$ws_ARM.Range("I88").Value = 2098.6365
$ws_ARM.Range("J88").Value = 1914.4445
$ws_ARM.Range("K88").Value = 2098.6365
$ws_ARM.Range("L88").Value = 1914.4445
$ws_ARM.Range("M88").Value = -1692.6365
$ws_ARM.Range("N88").Value = -2726.4445

$ws_ARM.Range("H91").Value = 2015.75
$ws_ARM.Range("I91").Value = 2098.6365
$ws_ARM.Range("J91").Value = 1914.4445
$ws_ARM.Range("K91").Value = 2098.6365
$ws_ARM.Range("L91").Value = 1914.4445
$ws_ARM.Range("M91").Value = -694.6365000000001
$ws_ARM.Range("N91").Value = -4722.4445

$ws_ARM.Range("H136").Value = 5862.4736
$ws_ARM.Range("I136").Value = 4360.3335
$ws_ARM.Range("J136").Value = 7214.4
$ws_ARM.Range("K136").Value = 13081.0005
$ws_ARM.Range("L136").Value = 21643.2
$ws_ARM.Range("M136").Value = -10531.0005
$ws_ARM.Range("N136").Value = -26743.2

$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_BSM.Range("H20").Value = 2635.3125
$ws_BSM.Range("I20").Value = 2816.9167
$ws_BSM.Range("J20").Value = 2090.5
$ws_BSM.Range("K20").Value = 2816.9167
$ws_BSM.Range("L20").Value = 2090.5
$ws_BSM.Range("M20").Value = -2569.9167
$ws_BSM.Range("N20").Value = -2584.5

$ws_BSM.Range("H99").Value = 31969.076
$ws_BSM.Range("I99").Value = 31969.076
$ws_BSM.Range("K99").Value = 31969.076
$ws_BSM.Range("M99").Value = -30471.076

$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CRP.Range("H31").Value = 5132.643
$ws_CRP.Range("I31").Value = 1629.6
$ws_CRP.Range("J31").Value = 7078.778
$ws_CRP.Range("K31").Value = 1629.6
$ws_CRP.Range("L31").Value = 7078.778
$ws_CRP.Range("M31").Value = -1334.6
$ws_CRP.Range("N31").Value = -7668.778

$ws_CRP.Range("H34").Value = 5132.643
$ws_CRP.Range("I34").Value = 1629.6
$ws_CRP.Range("J34").Value = 7078.778
$ws_CRP.Range("K34").Value = 1629.6
$ws_CRP.Range("L34").Value = 7078.778
$ws_CRP.Range("M34").Value = -1427.6
$ws_CRP.Range("N34").Value = -7482.778

$ws_CRP.Range("H122").Value = 50513184
$ws_CRP.Range("I122").Value = 84184904
$ws_CRP.Range("J122").Value = 5606.25
$ws_CRP.Range("K122").Value = 252554712
$ws_CRP.Range("L122").Value = 16818.75
$ws_CRP.Range("M122").Value = -252552262
$ws_CRP.Range("N122").Value = -21718.75

$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_CUL.Range("H4").Value = 58874772
$ws_CUL.Range("I4").Value = 71055550
$ws_CUL.Range("K4").Value = 213166650
$ws_CUL.Range("M4").Value = -213166538

$ws_CUL.Range("H29").Value = 814.2
$ws_CUL.Range("I29").Value = 1077.5714
$ws_CUL.Range("K29").Value = 3232.7142
$ws_CUL.Range("M29").Value = -2955.7142

$ws_CUL.Range("H34").Value = 494.45456
$ws_CUL.Range("I34").Value = 264.83334
$ws_CUL.Range("J34").Value = 770
$ws_CUL.Range("K34").Value = 794.5000200000001
$ws_CUL.Range("L34").Value = 2310
$ws_CUL.Range("M34").Value = -710.5000200000001
$ws_CUL.Range("N34").Value = -2478

$ws_CUL.Range("H63").Value = 148667.28
$ws_CUL.Range("I63").Value = 402503.8
$ws_CUL.Range("K63").Value = 1207511.4
$ws_CUL.Range("M63").Value = -1206762.4

$ws_CUL.Range("H66").Value = 148667.28
$ws_CUL.Range("I66").Value = 402503.8
$ws_CUL.Range("K66").Value = 3622534.2
$ws_CUL.Range("M66").Value = -3618790.2

$ws_CUL.Range("H81").Value = 0
$ws_CUL.Range("J81").Value = 0
$ws_CUL.Range("L81").ClearContents()
$ws_CUL.Range("N81").Value = 0

$ws_CUL.Range("H84").Value = 0
$ws_CUL.Range("J84").Value = 0
$ws_CUL.Range("L84").ClearContents()
$ws_CUL.Range("N84").Value = 0

$ws_CUL.Range("H103").Value = 4271.727
$ws_CUL.Range("I103").Value = 4497.25
$ws_CUL.Range("J103").Value = 4142.857
$ws_CUL.Range("K103").Value = 13491.75
$ws_CUL.Range("L103").Value = 12428.571
$ws_CUL.Range("M103").Value = -12612.75
$ws_CUL.Range("N103").Value = -14186.571

$ws_CUL.Range("H121").Value = 2026.129
$ws_CUL.Range("J121").Value = 2455.8
$ws_CUL.Range("L121").Value = 7367.400000000001
$ws_CUL.Range("N121").Value = -9987.400000000001

$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_GSM.Range("H69").Value = 50091
$ws_GSM.Range("I69").Value = 50182
$ws_GSM.Range("K69").Value = 50182
$ws_GSM.Range("M69").Value = -49433

$ws_GSM.Range("H72").Value = 50091
$ws_GSM.Range("I72").Value = 50182
$ws_GSM.Range("K72").Value = 150546
$ws_GSM.Range("M72").Value = -146802

$ws_GSM.Range("H109").Value = 49999.547
$ws_GSM.Range("J109").Value = 49999.547
$ws_GSM.Range("L109").Value = 49999.547
$ws_GSM.Range("N109").Value = -52079.547

$ws_GSM.Range("H122").Value = 3037.4
$ws_GSM.Range("I122").Value = 1778.2727
$ws_GSM.Range("K122").Value = 5334.8181
$ws_GSM.Range("M122").Value = -2884.8181

$ws_GSM.Range("H132").Value = 5953.05
$ws_GSM.Range("I132").Value = 5191.9
$ws_GSM.Range("J132").Value = 6714.2
$ws_GSM.Range("K132").Value = 15575.7
$ws_GSM.Range("L132").Value = 20142.6
$ws_GSM.Range("M132").Value = -13045.7
$ws_GSM.Range("N132").Value = -25202.6

$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_LTW.Range("H22").Value = 1028.7273
$ws_LTW.Range("I22").Value = 632.1111
$ws_LTW.Range("K22").Value = 632.1111
$ws_LTW.Range("M22").Value = -337.1111

$ws_LTW.Range("H27").Value = 1028.7273
$ws_LTW.Range("I27").Value = 632.1111
$ws_LTW.Range("K27").Value = 632.1111
$ws_LTW.Range("M27").Value = -525.1111

$ws_LTW.Range("H68").Value = 8406.200000000001
$ws_LTW.Range("I68").Value = 7725
$ws_LTW.Range("K68").Value = 7725
$ws_LTW.Range("M68").Value = -6976

$ws_LTW.Range("H71").Value = 8406.200000000001
$ws_LTW.Range("I71").Value = 7725
$ws_LTW.Range("K71").Value = 38625
$ws_LTW.Range("M71").Value = -34881

$ws_LTW.Range("H93").Value = 458267.1
$ws_LTW.Range("I93").Value = 3617.375
$ws_LTW.Range("J93").Value = 1670666.4
$ws_LTW.Range("K93").Value = 3617.375
$ws_LTW.Range("L93").Value = 1670666.4
$ws_LTW.Range("M93").Value = -2369.375
$ws_LTW.Range("N93").Value = -1673162.4

$ws_LTW.Range("H122").Value = 76927544
$ws_LTW.Range("I122").Value = 125002856
$ws_LTW.Range("K122").Value = 375008568
$ws_LTW.Range("M122").Value = -375006118

$ws_WVR = $wb.Worksheets.Item("WVR")
$ws_WVR.Range("H47").Value = 36183.26
$ws_WVR.Range("I47").Value = 36110.39
$ws_WVR.Range("J47").Value = 37495
$ws_WVR.Range("K47").Value = 36110.39
$ws_WVR.Range("L47").Value = 37495
$ws_WVR.Range("M47").Value = -35538.39
$ws_WVR.Range("N47").Value = -38639

$ws_WVR.Range("H81").Value = 14862
$ws_WVR.Range("I81").Value = 21979.2
$ws_WVR.Range("K81").Value = 43958.4
$ws_WVR.Range("M81").Value = -42897.4

$ws_WVR.Range("H84").Value = 14862
$ws_WVR.Range("I84").Value = 21979.2
$ws_WVR.Range("K84").Value = 219792
$ws_WVR.Range("M84").Value = -214488

$ws_WVR.Range("H126").Value = 3041.0833
$ws_WVR.Range("I126").Value = 2910.8
$ws_WVR.Range("K126").Value = 8732.400000000001
$ws_WVR.Range("M126").Value = -6262.400000000001
